# S1_Suivi_Modifications_DBA.xlsx
# "maj nomenclature et ajout des demandes de modif et etc."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Widen column D so the longer "Description" text fits -----------------
$ws.Columns.Item(4).ColumnWidth = 35.25

# --- Row 8: first change-request entry -------------------------------------
$ws.Range("A8").Value = "DA1"
$ws.Range("B8").Value = "David Paquet"
$ws.Range("C8").Value = 42257
$ws.Range("C8").NumberFormat = "mm-dd-yy"
$ws.Range("D8").Value = "CodeSysExp de varchar(10) à varchar(20)"
$ws.Range("E8").Value = "Annulée"

# --- Row 9: second change-request entry ------------------------------------
$ws.Range("A9").Value = "GA01"
$ws.Range("B9").Value = "Gabriel Simard"
$ws.Range("C9").Value = 42257
# Reuse C8's format (copy/paste-special) so both date cells share one style
# entry instead of the engine minting a duplicate cellXfs record.
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D9").Value = "DescJeu de varchar(250) à varchar(350)"
$ws.Range("E9").Value = "Terminée"

# --- Move the active selection to match the saved view ---------------------
$ws.Range("K23").Select()
